# Natmi following Dr Hou advice
# Update recomputed Cxcl12-Itgb3 LR-pair statistics (rows 2-10, columns E:T)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; F=1; G=130.955829;        H=392.867487;        I=0.5336535908353144; J=0.5336535908353144; K=3; L=1; M=8.970048;           N=26.910144;          O=0.487108783009476;  P=0.4871087830094759; Q=1174.680072009792; R=10572.12064808813;  S=0.2599473511804268;  T=0.2599473511804268  }
    3  = @{ E=3; F=1; G=130.955829;        H=392.867487;        I=0.5336535908353144; J=0.5336535908353144; K=3; L=1; M=9.012070666666666;  N=27.036212;          O=0.489390778604016;  P=0.489390778604016;  Q=1180.183185159916; R=10621.64866643924;  S=0.2611651463237235;  T=0.2611651463237235  }
    4  = @{ E=3; F=1; G=130.955829;        H=392.867487;        I=0.5336535908353144; J=0.5336535908353144; K=3; L=1; M=0.4327576666666667; N=1.298273;           O=0.02350043838650813;P=0.02350043838650813;Q=56.672138994439;    R=510.049250949951;   S=0.01254109333116412; T=0.01254109333116412 }
    5  = @{ E=3; F=1; G=66.39541;          H=199.18623;         I=0.2705656497465488; J=0.2705656497465488; K=3; L=1; M=8.970048;           N=26.910144;          O=0.487108783009476;  P=0.4871087830094759; Q=595.57001467968;    R=5360.13013211712;   S=0.1317949043722095;  T=0.1317949043722095  }
    6  = @{ E=3; F=1; G=66.39541;          H=199.18623;         I=0.2705656497465488; J=0.2705656497465488; K=3; L=1; M=9.012070666666666;  N=27.036212;          O=0.489390778604016;  P=0.489390778604016;  Q=598.3601268623066;  R=5385.24114176076;   S=0.132412333992965;   T=0.132412333992965   }
    7  = @{ E=3; F=1; G=66.39541;          H=199.18623;         I=0.2705656497465488; J=0.2705656497465488; K=3; L=1; M=0.4327576666666667; N=1.298273;           O=0.02350043838650813;P=0.02350043838650813;Q=28.73312270897667;  R=258.59810438079;    S=0.006358411381374309;T=0.006358411381374309}
    8  = @{ E=3; F=1; G=48.043585;         H=144.130755;        I=0.1957807594181367; J=0.1957807594181367; K=3; L=1; M=8.970048;           N=26.910144;          O=0.487108783009476;  P=0.4871087830094759; Q=430.95326354208;    R=3878.579371878719;  S=0.09536652745683957; T=0.09536652745683956 }
    9  = @{ E=3; F=1; G=48.043585;         H=144.130755;        I=0.1957807594181367; J=0.1957807594181367; K=3; L=1; M=9.012070666666666;  N=27.036212;          O=0.489390778604016;  P=0.489390778604016;  Q=432.9721831000066;  R=3896.74964790006;   S=0.09581329828732746; T=0.09581329828732744 }
    10 = @{ E=3; F=1; G=48.043585;         H=144.130755;        I=0.1957807594181367; J=0.1957807594181367; K=3; L=1; M=0.4327576666666667; N=1.298273;           O=0.02350043838650813;P=0.02350043838650813;Q=20.79122974290167;  R=187.121067686115;   S=0.004600933673969692;T=0.004600933673969691}
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
